{"js": "// Replace the old \"Datas das campanhas de 2018 que usam Perseu: 30 de outubro\n// a 8 de novembro e 29 de novembro a 8 de dezembro\" paragraphs with the new,\n// translated \"Datas das campanhas de Taurus: 16-25 de janeiro\" text. The old\n// paragraphs are built from four separate runs; the new paragraph is a single\n// plain run (no run formatting carried over), so each matching paragraph is\n// cleared first and then given fresh, unformatted text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst oldLead = \"Datas das campanhas de 2018 que usam\";\nconst newText = \"Datas das campanhas de Taurus: 16-25 de janeiro\";\n\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(oldLead) !== -1) {\n    targets.push(p);\n  }\n}\n\n// Clear each matching paragraph's content (drops the old runs/formatting).\nfor (const p of targets) {\n  p.clear();\n}\nawait context.sync();\n\n// Insert the new plain-text run into each cleared paragraph.\nfor (const p of targets) {\n  p.insertText(newText, \"Start\");\n}\nawait context.sync();\n", "ps1": "# Replace the old \"Datas das campanhas de 2018 que usam Perseu: 30 de outubro\n# a 8 de novembro e 29 de novembro a 8 de dezembro\" paragraphs with the new,\n# translated \"Datas das campanhas de Taurus: 16-25 de janeiro\" text. The old\n# paragraphs are built from four separate runs; the new paragraph is a single\n# plain run (no leftover run formatting), so each matching paragraph's range\n# (excluding its trailing paragraph mark) is deleted outright and then given\n# fresh, unformatted text via InsertAfter.\n\n$d = $word.ActiveDocument\n$oldLead = \"Datas das campanhas de 2018 que usam\"\n$newText = \"Datas das campanhas de Taurus: 16-25 de janeiro\"\n\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$oldLead*\") {\n        $targets += $p\n    }\n}\n\nforeach ($p in $targets) {\n    $r = $p.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Delete() | Out-Null\n    $r.InsertAfter($newText) | Out-Null\n}\n"}
